$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Artha "
$ws.Range("B3").Value = "Asih "
$ws.Range("B4").Value = "Bantar "
$ws.Range("B5").Value = "Berdikari"

$ws.Range("B6").Select()
